# Auto-generated: apply scheduled price-refresh updates to Typhon_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4771987.5
$ws.Range("J17").Value = 5274165.5
$ws.Range("L17").Value = 15822496.5
$ws.Range("N17").Value = -15822832.5
$ws.Range("H19").Value = 2132.9167
$ws.Range("I19").Value = 4552.2
$ws.Range("J19").Value = 404.85715
$ws.Range("K19").Value = 4552.2
$ws.Range("L19").Value = 404.85715
$ws.Range("M19").Value = -4377.2
$ws.Range("N19").Value = -754.85715
$ws.Range("H116").Value = 4334.278
$ws.Range("I116").Value = 1938.125
$ws.Range("J116").Value = 6251.2
$ws.Range("K116").Value = 1938.125
$ws.Range("L116").Value = 6251.2
$ws.Range("M116").Value = 1503.875
$ws.Range("N116").Value = -13135.2
$ws.Range("H129").Value = 286939.9
$ws.Range("J129").Value = 346232.3
$ws.Range("L129").Value = 1038696.9
$ws.Range("N129").Value = -1048696.9
$ws.Range("H135").Value = 11114307
$ws.Range("I135").Value = 552.2973
$ws.Range("J135").Value = 62515424
$ws.Range("K135").Value = 4970.6757
$ws.Range("L135").Value = 562638816
$ws.Range("M135").Value = -2435.6757
$ws.Range("N135").Value = -562643886
$ws.Range("H137").Value = 1638.279
$ws.Range("I137").Value = 1777.3448
$ws.Range("J137").Value = 1350.2142
$ws.Range("K137").Value = 5332.0344
$ws.Range("L137").Value = 4050.6426
$ws.Range("M137").Value = -2782.0344
$ws.Range("N137").Value = -9150.642599999999
$ws.Range("H138").Value = 22224912
$ws.Range("I138").Value = 41668150
$ws.Range("J138").Value = 4065.524
$ws.Range("K138").Value = 125004450
$ws.Range("L138").Value = 12196.572
$ws.Range("M138").Value = -124999310
$ws.Range("N138").Value = -22476.572
$ws.Range("H141").Value = 2094.4482
$ws.Range("I141").Value = 1237
$ws.Range("J141").Value = 3999.889
$ws.Range("K141").Value = 3711
$ws.Range("L141").Value = 11999.667
$ws.Range("M141").Value = 1469
$ws.Range("N141").Value = -22359.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1092.2693
$ws.Range("I2").Value = 1187.2354
$ws.Range("J2").Value = 912.8889
$ws.Range("K2").Value = 1187.2354
$ws.Range("L2").Value = 912.8889
$ws.Range("M2").Value = -1074.2354
$ws.Range("N2").Value = -1138.8889
$ws.Range("H45").Value = 2465.353
$ws.Range("I45").Value = 2348.7144
$ws.Range("J45").Value = 2653.7693
$ws.Range("K45").Value = 2348.7144
$ws.Range("L45").Value = 2653.7693
$ws.Range("M45").Value = -1971.7144
$ws.Range("N45").Value = -3407.7693
$ws.Range("H61").Value = 347162.28
$ws.Range("I61").Value = 375932.9
$ws.Range("J61").Value = 1914.75
$ws.Range("K61").Value = 375932.9
$ws.Range("L61").Value = 1914.75
$ws.Range("M61").Value = -375720.9
$ws.Range("N61").Value = -2338.75
$ws.Range("H74").Value = 37040016
$ws.Range("I74").Value = 50002796
$ws.Range("J74").Value = 3499.7144
$ws.Range("K74").Value = 50002796
$ws.Range("L74").Value = 3499.7144
$ws.Range("M74").Value = -50001922
$ws.Range("N74").Value = -5247.7144
$ws.Range("H77").Value = 37040016
$ws.Range("I77").Value = 50002796
$ws.Range("J77").Value = 3499.7144
$ws.Range("K77").Value = 250013980
$ws.Range("L77").Value = 17498.572
$ws.Range("M77").Value = -250009612
$ws.Range("N77").Value = -26234.572
$ws.Range("H116").Value = 1092.2693
$ws.Range("I116").Value = 1187.2354
$ws.Range("J116").Value = 912.8889
$ws.Range("K116").Value = 1187.2354
$ws.Range("L116").Value = 912.8889
$ws.Range("M116").Value = 1106.7646
$ws.Range("N116").Value = -5500.8889
$ws.Range("H136").Value = 347162.28
$ws.Range("I136").Value = 375932.9
$ws.Range("J136").Value = 1914.75
$ws.Range("K136").Value = 1127798.7
$ws.Range("L136").Value = 5744.25
$ws.Range("M136").Value = -1125248.7
$ws.Range("N136").Value = -10844.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1092.2693
$ws.Range("I3").Value = 1187.2354
$ws.Range("J3").Value = 912.8889
$ws.Range("K3").Value = 1187.2354
$ws.Range("L3").Value = 912.8889
$ws.Range("M3").Value = -1073.2354
$ws.Range("N3").Value = -1140.8889
$ws.Range("H20").Value = 1986
$ws.Range("I20").Value = 2243.25
$ws.Range("K20").Value = 2243.25
$ws.Range("M20").Value = -1996.25
$ws.Range("H134").Value = 2221.4814
$ws.Range("I134").Value = 2341.6326
$ws.Range("J134").Value = 1044
$ws.Range("K134").Value = 7024.8978
$ws.Range("L134").Value = 3132
$ws.Range("M134").Value = -4489.8978
$ws.Range("N134").Value = -8202

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4831.5757
$ws.Range("I31").Value = 4073.8572
$ws.Range("K31").Value = 4073.8572
$ws.Range("M31").Value = -3778.8572
$ws.Range("H34").Value = 4831.5757
$ws.Range("I34").Value = 4073.8572
$ws.Range("K34").Value = 4073.8572
$ws.Range("M34").Value = -3871.8572
$ws.Range("H58").Value = 11066.02
$ws.Range("I58").Value = 903.25
$ws.Range("J58").Value = 37198.855
$ws.Range("K58").Value = 903.25
$ws.Range("L58").Value = 37198.855
$ws.Range("M58").Value = -700.25
$ws.Range("N58").Value = -37604.855
$ws.Range("H132").Value = 1659.0807
$ws.Range("I132").Value = 1167.7778
$ws.Range("J132").Value = 4975.375
$ws.Range("K132").Value = 3503.3334
$ws.Range("L132").Value = 14926.125
$ws.Range("M132").Value = -973.3334000000004
$ws.Range("N132").Value = -19986.125
$ws.Range("H134").Value = 746.4400000000001
$ws.Range("I134").Value = 604.9524
$ws.Range("K134").Value = 1814.8572
$ws.Range("M134").Value = 720.1428000000001
$ws.Range("H136").Value = 11066.02
$ws.Range("I136").Value = 903.25
$ws.Range("J136").Value = 37198.855
$ws.Range("K136").Value = 2709.75
$ws.Range("L136").Value = 111596.565
$ws.Range("M136").Value = -159.75
$ws.Range("N136").Value = -116696.565
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 4141.636
$ws.Range("J106").Value = 4141.636
$ws.Range("L106").Value = 12424.908
$ws.Range("N106").Value = -14316.908
$ws.Range("H131").Value = 122745.34
$ws.Range("J131").Value = 132377.88
$ws.Range("L131").Value = 397133.64
$ws.Range("N131").Value = -407213.64

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 184.74074
$ws.Range("J55").Value = 247.25
$ws.Range("L55").Value = 247.25
$ws.Range("N55").Value = -593.25
$ws.Range("H122").Value = 579055.0600000001
$ws.Range("I122").Value = 1034162.3
$ws.Range("J122").Value = 2585.9333
$ws.Range("K122").Value = 3102486.9
$ws.Range("L122").Value = 7757.7999
$ws.Range("M122").Value = -3100036.9
$ws.Range("N122").Value = -12657.7999
$ws.Range("H132").Value = 2352.72
$ws.Range("J132").Value = 5416.3335
$ws.Range("L132").Value = 16249.0005
$ws.Range("N132").Value = -21309.0005
$ws.Range("H136").Value = 1159.4584
$ws.Range("I136").Value = 1166.3914
$ws.Range("K136").Value = 3499.1742
$ws.Range("M136").Value = -949.1741999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1165.8975
$ws.Range("I132").Value = 789.1
$ws.Range("J132").Value = 2421.889
$ws.Range("K132").Value = 2367.3
$ws.Range("L132").Value = 7265.667
$ws.Range("M132").Value = 162.6999999999998
$ws.Range("N132").Value = -12325.667
$ws.Range("H136").Value = 22729994
$ws.Range("I136").Value = 34484068
$ws.Range("J136").Value = 5447.3335
$ws.Range("K136").Value = 103452204
$ws.Range("L136").Value = 16342.0005
$ws.Range("M136").Value = -103449654
$ws.Range("N136").Value = -21442.0005

Write-Output "All changes applied"
